$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from H1 into the two new
# header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-20: new I0 / IF columns.
$data = @(
    @(5,5),
    @(7,8),
    @(7,7),
    @(6,6),
    @(7,7),
    @(5,5),
    @(7,7),
    @(9,9),
    @(3,4),
    @(8,8),
    @(5,5),
    @(6,6),
    @(6,7),
    @(6,6),
    @(7,7),
    @(9,9),
    @(9,9),
    @(6,6),
    @(9,9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
